# Fixed a bug in stats2
# This applies a row permutation (rows 2-21) to the "symbol"/"reelN" table
# on the active worksheet, matching the corrected stats2 output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(401, 9, 48, 67, 75, 45)
    3  = @(601, 9, 60, 67, 60, 42)
    4  = @(801, 3, 67, 65, 52, 45)
    5  = @(101, 9, 30, 15, 60, 15)
    6  = @(901, 16, 15, 45, 60, 60)
    7  = @(1001, 18, 30, 75, 60, 72)
    8  = @(1203, 3, 15, 15, 15, 15)
    9  = @(902, 1, 0, 0, 0, 0)
    10 = @(301, 6, 45, 30, 60, 45)
    11 = @(1202, 2, 10, 10, 10, 10)
    12 = @(1201, 2, 10, 10, 10, 10)
    13 = @(501, 9, 52, 30, 75, 45)
    14 = @(701, 3, 90, 45, 97, 15)
    15 = @(201, 9, 30, 15, 45, 30)
    16 = @(1, 0, 2, 2, 2, 2)
    17 = @(3, 0, 3, 3, 3, 3)
    18 = @(502, 0, 4, 0, 0, 0)
    19 = @(2, 0, 2, 2, 2, 2)
    20 = @(1101, 0, 15, 30, 30, 0)
    21 = @(802, 0, 4, 5, 4, 0)
}

foreach ($rowNum in $data.Keys) {
    $values = $data[$rowNum]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $values[$col - 1]
    }
}
